$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 74 with revised figures
$ws.Range("B74").Value = 57157
$ws.Range("C74").Value = 4749
$ws.Range("D74").Value = 71
$ws.Range("E74").Value = 4679
$ws.Range("F74").Value = 52408
$ws.Range("G74").Value = 4165
$ws.Range("H74").Value = 48243

# Add new row 75 for quarter 01-04-2021
$ws.Range("A75").NumberFormat = "@"
$ws.Range("A75").Value = "01-04-2021"
$ws.Range("A75").Style = "Normal"
$ws.Range("B75").Value = 56685
$ws.Range("C75").Value = 4769
$ws.Range("D75").Value = 71
$ws.Range("E75").Value = 4698
$ws.Range("F75").Value = 51916
$ws.Range("G75").Value = 3497
$ws.Range("H75").Value = 48419
